$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.730.60"
$ws.Range("E2").Value = "'  +4.66%  "
$ws.Range("D3").Value = "'3.057.83"
$ws.Range("E3").Value = "'  +2.37%  "
$ws.Range("E4").Value = "'  +0.31%  "
$ws.Range("D5").Value = "'575.90"
$ws.Range("E5").Value = "'  +2.53%  "
$ws.Range("D6").Value = "'142.35"
$ws.Range("E6").Value = "'  +3.66%  "
$ws.Range("E7").Value = "'  +0.07%  "
$ws.Range("D8").Value = "'3.046.43"
$ws.Range("E8").Value = "'  +2.28%  "
$ws.Range("E9").Value = "'  +1.09%  "
$ws.Range("D10").Value = "'0.138"
$ws.Range("E10").Value = "'  +4.96%  "
$ws.Range("D11").Value = "'5.46"
$ws.Range("E11").Value = "'  +12.50%  "
$ws.Range("E12").Value = "'  +1.64%  "
$ws.Range("D13").Value = "'0.0000239"
$ws.Range("E13").Value = "'  +4.83%  "
$ws.Range("D14").Value = "'34.70"
$ws.Range("E14").Value = "'  +3.35%  "
$ws.Range("E15").Value = "'  -0.19%  "
$ws.Range("D16").Value = "'3.567.22"
$ws.Range("E16").Value = "'  +2.94%  "
$ws.Range("E17").Value = "'  +2.86%  "
$ws.Range("D18").Value = "'3.054.47"
$ws.Range("E18").Value = "'  +2.61%  "
$ws.Range("D19").Value = "'61.800.86"
$ws.Range("E19").Value = "'  +5.09%  "
$ws.Range("D20").Value = "'448.73"
$ws.Range("E20").Value = "'  +6.16%  "
$ws.Range("D21").Value = "'13.91"
$ws.Range("E21").Value = "'  +2.85%  "
$ws.Range("E22").Value = "'  +2.47%  "
$ws.Range("D23").Value = "'7.26"
$ws.Range("E23").Value = "'  +1.96%  "
$ws.Range("D24").Value = "'13.64"
$ws.Range("E24").Value = "'  +1.40%  "
$ws.Range("D25").Value = "'81.81"
$ws.Range("E25").Value = "'  +1.84%  "
$ws.Range("E26").Value = "'  +0.15%  "
$ws.Range("E27").Value = "'  +5.36%  "
$ws.Range("E28").Value = "'  +0.29%  "
$ws.Range("E29").Value = "'  +3.99%  "
$ws.Range("D30").Value = "'8.01"
$ws.Range("E30").Value = "'  +3.33%  "
$ws.Range("D31").Value = "'6.59"
$ws.Range("E31").Value = "'  +7.92%  "
$ws.Range("D32").Value = "'26.58"
$ws.Range("E32").Value = "'  +3.87%  "
$ws.Range("E33").Value = "'  +7.09%  "
$ws.Range("D34").Value = "'0.0₃0805"
$ws.Range("E34").Value = "'  +4.86%  "
$ws.Range("E35").Value = "'  +2.32%  "
$ws.Range("D36").Value = "'6.06"
$ws.Range("E36").Value = "'  +5.79%  "
$ws.Range("E37").Value = "'  +5.18%  "
$ws.Range("D38").Value = "'50.13"
$ws.Range("E38").Value = "'  +2.53%  "
$ws.Range("D39").Value = "'2.94"
$ws.Range("E39").Value = "'  +6.32%  "
$ws.Range("D40").Value = "'8.82"
$ws.Range("E40").Value = "'  +1.90%  "
$ws.Range("D41").Value = "'413.34"
$ws.Range("E41").Value = "'  +3.21%  "
$ws.Range("D42").Value = "'0.0366"
$ws.Range("E42").Value = "'  +5.19%  "
$ws.Range("D43").Value = "'2.767.21"
$ws.Range("E43").Value = "'  +0.87%  "
$ws.Range("E44").Value = "'  +0.49%  "
$ws.Range("D45").Value = "'0.263"
$ws.Range("E45").Value = "'  +7.25%  "
$ws.Range("D46").Value = "'36.60"
$ws.Range("E46").Value = "'  +13.22%  "
$ws.Range("E48").Value = "'  +3.19%  "
$ws.Range("D49").Value = "'123.04"
$ws.Range("E49").Value = "'  -1.56%  "
$ws.Range("E50").Value = "'  +1.42%  "
$ws.Range("D51").Value = "'23.98"
$ws.Range("E51").Value = "'  +2.89%  "

# Reset number format/style so text-coercion via leading apostrophe
# does not leave a residual "Text" style on the cells (keep original formatting).
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Style = "Normal"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Style = "Normal"
